$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE (column N) for rows 2, 3, 4
$ws.Range("N2").Value = 51.05762969290213
$ws.Range("N3").Value = 51.05762969290213
$ws.Range("N4").Value = 51.05762969290213

# Update Gold Feb 26 (row 4) price/RSI/5-day return figures
$ws.Range("D4").Value = 4254.1
$ws.Range("E4").Value = 72.7
$ws.Range("F4").Value = 0.85
